$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Comment/SamplePortion/SamplePortionUnit columns (J, K, L) are being
# reordered so that SamplePortion and SamplePortionUnit come before
# Comment (Comment moves to the end), for every row in the header block
# (machine name, french label, type, format description, example).
# This is a simple left-rotation of the J/K/L column values in each row:
#   new J = old K
#   new K = old L
#   new L = old J
1..5 | ForEach-Object {
    $row = $_
    $oldJ = $ws.Cells.Item($row, 10).Value()
    $oldK = $ws.Cells.Item($row, 11).Value()
    $oldL = $ws.Cells.Item($row, 12).Value()

    $ws.Cells.Item($row, 10).Value = $oldK
    $ws.Cells.Item($row, 11).Value = $oldL
    $ws.Cells.Item($row, 12).Value = $oldJ
}
